$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update logger / thread name labels
$ws.Range("A2").Value = "testlogger"
$ws.Range("A3").Value = "testing"

# Swap DEBUG / FATAL labels
$ws.Range("A6").Value = "FATAL"
$ws.Range("A7").Value = "DEBUG"

# Update counts
$ws.Range("B2").Value = 30.0
$ws.Range("B3").Value = 30.0
$ws.Range("B4").Value = 6.0
$ws.Range("B5").Value = 6.0
$ws.Range("B6").Value = 6.0
$ws.Range("B7").Value = 6.0
$ws.Range("B8").Value = 6.0
